$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Volumen/Precio" block (columns D,H,I,J,K,L,M,P) for rows 2,3,6,7,8,9,10,11
# has been rotated through a single cycle: 2 -> 3 -> 11 -> 8 -> 10 -> 7 -> 6 -> 9 -> 2
# i.e. each row now carries the data that used to belong to the "next" row in that
# cycle. We capture the original values first, then write the rotated values back.

$rows = @(2, 3, 11, 8, 10, 7, 6, 9)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        H = $ws.Cells.Item($r, 8).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

for ($i = 0; $i -lt $rows.Length; $i++) {
    $dest = $rows[$i]
    $src = $rows[($i + 1) % $rows.Length]
    $vals = $orig[$src]

    $ws.Cells.Item($dest, 4).Value2 = $vals.D
    $ws.Cells.Item($dest, 8).Value2 = $vals.H
    $ws.Cells.Item($dest, 9).Value2 = $vals.I
    $ws.Cells.Item($dest, 10).Value2 = $vals.J
    $ws.Cells.Item($dest, 11).Value2 = $vals.K
    $ws.Cells.Item($dest, 12).Value2 = $vals.L
    $ws.Cells.Item($dest, 13).Value2 = $vals.M
    $ws.Cells.Item($dest, 16).Value2 = $vals.P
}
